# Auto-generated: applies 2022-08-03 daily crime-count update (column I = year 2022 YTD totals)
# to the Chicago violent-crime workbook. For each affected worksheet, sets the listed
# cells in column I to their new running totals, per the commit diff.

$wb = $excel.ActiveWorkbook

$updates = @{
    'Citywide Totals' = @{ 'I2' = 4132; 'I3' = 4302; 'I4' = 1001; 'I5' = 394; 'I6' = 4749; 'I7' = 14578 }
    'West Ridge' = @{ 'I3' = 44; 'I7' = 155 }
    'Bridgeport' = @{ 'I2' = 27; 'I7' = 79 }
    'Fuller Park' = @{ 'I2' = 13; 'I3' = 17 }
    'Grand Crossing' = @{ 'I2' = 149; 'I7' = 471 }
    'North Lawndale' = @{ 'I2' = 133; 'I6' = 185; 'I7' = 565 }
    'Gage Park' = @{ 'I6' = 49; 'I7' = 141 }
    'South Deering' = @{ 'I6' = 32; 'I7' = 125 }
    'New City' = @{ 'I2' = 106; 'I3' = 94; 'I6' = 99; 'I7' = 329 }
    'By Neighborhood' = @{ 'I3' = 10; 'I7' = 462; 'I8' = 878; 'I14' = 79; 'I18' = 103; 'I19' = 403; 'I20' = 356; 'I21' = 80; 'I25' = 71; 'I27' = 132; 'I29' = 930; 'I31' = 141; 'I33' = 671; 'I36' = 201; 'I37' = 471; 'I41' = 65; 'I42' = 499; 'I47' = 98; 'I50' = 63; 'I51' = 151; 'I52' = 318; 'I55' = 164; 'I63' = 56; 'I65' = 329; 'I67' = 565; 'I74' = 29; 'I76' = 214; 'I78' = 206; 'I79' = 395; 'I83' = 298; 'I84' = 125; 'I85' = 660; 'I95' = 240; 'I96' = 155; 'I97' = 110; 'I101' = 14578 }
    'South Chicago' = @{ 'I2' = 103; 'I7' = 298 }
    'West Pullman' = @{ 'I3' = 92; 'I6' = 41; 'I7' = 240 }
    'Garfield Park' = @{ 'I5' = 27; 'I6' = 207; 'I7' = 671 }
    'Englewood' = @{ 'I2' = 271; 'I3' = 320; 'I4' = 47; 'I6' = 255; 'I7' = 930 }
    'Chatham' = @{ 'I2' = 147; 'I6' = 111; 'I7' = 403 }
    'River North' = @{ 'I6' = 89; 'I7' = 214 }
    'South Shore' = @{ 'I2' = 170; 'I3' = 267; 'I5' = 23; 'I7' = 660 }
    'Hermosa' = @{ 'I2' = 21; 'I7' = 65 }
    'Humboldt Park' = @{ 'I3' = 172; 'I7' = 499 }
    'Rogers Park' = @{ 'I3' = 53; 'I6' = 81; 'I7' = 206 }
    'Lower West Side' = @{ 'I3' = 50; 'I4' = 9; 'I7' = 164 }
    'Chinatown' = @{ 'I3' = 12; 'I7' = 80 }
    'Roseland' = @{ 'I2' = 117; 'I4' = 26; 'I6' = 116; 'I7' = 395 }
    'Chicago Lawn' = @{ 'I3' = 111; 'I7' = 356 }
    'Calumet Heights' = @{ 'I6' = 45; 'I7' = 103 }
    'Grand Boulevard' = @{ 'I2' = 66; 'I7' = 201 }
    'Little Village' = @{ 'I2' = 90; 'I7' = 318 }
    'East Side' = @{ 'I3' = 19; 'I7' = 71 }
    'Kenwood' = @{ 'I6' = 36; 'I7' = 98 }
    'Lincoln Square' = @{ 'I3' = 16; 'I7' = 63 }
    'West Town' = @{ 'I6' = 64; 'I7' = 110 }
    'Austin' = @{ 'I2' = 278; 'I3' = 248; 'I6' = 278; 'I7' = 878 }
    'Edgewater' = @{ 'I6' = 53; 'I7' = 132 }
    'Little Italy, UIC' = @{ 'I6' = 60; 'I7' = 151 }
    'Andersonville' = @{ 'I2' = 4; 'I6' = 10 }
    'Auburn Gresham' = @{ 'I2' = 158; 'I3' = 144; 'I7' = 462 }
    'Printers Row' = @{ 'I6' = 20; 'I7' = 29 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellValues = $updates[$sheetName]
    foreach ($cellRef in $cellValues.Keys) {
        $ws.Range($cellRef).Value = $cellValues[$cellRef]
    }
}

Write-Host "Applied $($updates.Values | ForEach-Object { $_.Keys.Count } | Measure-Object -Sum | Select-Object -ExpandProperty Sum) cell updates across $($updates.Keys.Count) worksheets."
